$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix corrupted commas in contractor/company name entries ---
# (scraping artifact replaced stray commas with periods in a handful of
#  "Razon social" / "Nombre Fantasia" cells)
$ws.Range("E34").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F34").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E42").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F42").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E71").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("F71").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("E77").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F77").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E79").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E92").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F92").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E106").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E122").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E152").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E153").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F153").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E155").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# --- Fix floating point number formatting in the "Importe" column (H) ---
# Values were stored as localized (es-AR) formatted text, e.g. "2.073,00"
# (dot thousands separator, comma decimal separator). Re-write them as plain
# decimal text, e.g. "2073.00", without letting Excel coerce them to numbers.
$importeRange = $ws.Range("H2:H189")
$importeRange.NumberFormat = "@"
$ws.Range("H2").Value = "2073.00"
$ws.Range("H3").Value = "340.00"
$ws.Range("H4").Value = "682.00"
$ws.Range("H5").Value = "1180.00"
$ws.Range("H6").Value = "32.00"
$ws.Range("H7").Value = "17060.00"
$ws.Range("H8").Value = "33755.00"
$ws.Range("H9").Value = "11453.00"
$ws.Range("H10").Value = "11550.00"
$ws.Range("H11").Value = "1770.00"
$ws.Range("H12").Value = "748.00"
$ws.Range("H13").Value = "43.80"
$ws.Range("H14").Value = "110.00"
$ws.Range("H15").Value = "29053.37"
$ws.Range("H16").Value = "39838.12"
$ws.Range("H17").Value = "63495.81"
$ws.Range("H18").Value = "11953.85"
$ws.Range("H19").Value = "7828.39"
$ws.Range("H20").Value = "10787.54"
$ws.Range("H21").Value = "7425.00"
$ws.Range("H22").Value = "1034.40"
$ws.Range("H23").Value = "25882.39"
$ws.Range("H24").Value = "3968.81"
$ws.Range("H25").Value = "164.00"
$ws.Range("H26").Value = "234.00"
$ws.Range("H27").Value = "248.46"
$ws.Range("H28").Value = "6560.98"
$ws.Range("H29").Value = "390.00"
$ws.Range("H30").Value = "618.00"
$ws.Range("H31").Value = "790.00"
$ws.Range("H32").Value = "1153.90"
$ws.Range("H33").Value = "1003.88"
$ws.Range("H34").Value = "270.20"
$ws.Range("H35").Value = "1059.97"
$ws.Range("H36").Value = "80.67"
$ws.Range("H37").Value = "50.00"
$ws.Range("H38").Value = "6342.50"
$ws.Range("H39").Value = "0.12"
$ws.Range("H40").Value = "210.00"
$ws.Range("H41").Value = "65.58"
$ws.Range("H42").Value = "54.62"
$ws.Range("H43").Value = "554.98"
$ws.Range("H44").Value = "1930.00"
$ws.Range("H45").Value = "40051.75"
$ws.Range("H46").Value = "508.20"
$ws.Range("H47").Value = "4378.90"
$ws.Range("H48").Value = "3753.00"
$ws.Range("H49").Value = "566.39"
$ws.Range("H50").Value = "19850.69"
$ws.Range("H51").Value = "2583.00"
$ws.Range("H52").Value = "36.48"
$ws.Range("H53").Value = "20.00"
$ws.Range("H54").Value = "799.50"
$ws.Range("H55").Value = "130.00"
$ws.Range("H56").Value = "109.00"
$ws.Range("H57").Value = "123.68"
$ws.Range("H58").Value = "101.02"
$ws.Range("H59").Value = "193.88"
$ws.Range("H60").Value = "1720.21"
$ws.Range("H61").Value = "933.74"
$ws.Range("H62").Value = "522.00"
$ws.Range("H63").Value = "4163.50"
$ws.Range("H64").Value = "224.00"
$ws.Range("H65").Value = "609.34"
$ws.Range("H66").Value = "13035.50"
$ws.Range("H67").Value = "8800.00"
$ws.Range("H68").Value = "390.00"
$ws.Range("H69").Value = "6621.80"
$ws.Range("H70").Value = "330.33"
$ws.Range("H71").Value = "587.00"
$ws.Range("H72").Value = "590.00"
$ws.Range("H73").Value = "432.00"
$ws.Range("H74").Value = "10449.00"
$ws.Range("H75").Value = "11050.00"
$ws.Range("H76").Value = "3428.80"
$ws.Range("H77").Value = "2644.08"
$ws.Range("H78").Value = "260.00"
$ws.Range("H79").Value = "854.00"
$ws.Range("H80").Value = "363.00"
$ws.Range("H81").Value = "393.10"
$ws.Range("H82").Value = "250.00"
$ws.Range("H83").Value = "383.00"
$ws.Range("H84").Value = "225368.00"
$ws.Range("H85").Value = "171.30"
$ws.Range("H86").Value = "16000.00"
$ws.Range("H87").Value = "3150.00"
$ws.Range("H88").Value = "56.38"
$ws.Range("H89").Value = "14.02"
$ws.Range("H90").Value = "218.93"
$ws.Range("H91").Value = "120.00"
$ws.Range("H92").Value = "76.53"
$ws.Range("H93").Value = "321.24"
$ws.Range("H94").Value = "10.26"
$ws.Range("H95").Value = "2427.29"
$ws.Range("H96").Value = "50.00"
$ws.Range("H97").Value = "379.08"
$ws.Range("H98").Value = "796.50"
$ws.Range("H99").Value = "6.38"
$ws.Range("H100").Value = "792.00"
$ws.Range("H101").Value = "21.25"
$ws.Range("H102").Value = "5969.00"
$ws.Range("H103").Value = "4559.30"
$ws.Range("H104").Value = "435.60"
$ws.Range("H105").Value = "1006.20"
$ws.Range("H106").Value = "5524.50"
$ws.Range("H107").Value = "352.80"
$ws.Range("H108").Value = "5052.00"
$ws.Range("H109").Value = "8326.56"
$ws.Range("H110").Value = "1890.00"
$ws.Range("H111").Value = "249.00"
$ws.Range("H112").Value = "167.92"
$ws.Range("H113").Value = "294.00"
$ws.Range("H114").Value = "440.00"
$ws.Range("H115").Value = "270.00"
$ws.Range("H116").Value = "150.00"
$ws.Range("H117").Value = "173.00"
$ws.Range("H118").Value = "300.00"
$ws.Range("H119").Value = "900.00"
$ws.Range("H120").Value = "100.00"
$ws.Range("H121").Value = "130.00"
$ws.Range("H122").Value = "12900.00"
$ws.Range("H123").Value = "270.00"
$ws.Range("H124").Value = "250.00"
$ws.Range("H125").Value = "500.00"
$ws.Range("H126").Value = "98.28"
$ws.Range("H127").Value = "339.00"
$ws.Range("H128").Value = "910.96"
$ws.Range("H129").Value = "1302.78"
$ws.Range("H130").Value = "779594.20"
$ws.Range("H131").Value = "759000.00"
$ws.Range("H132").Value = "3276.40"
$ws.Range("H133").Value = "5008.74"
$ws.Range("H134").Value = "200.00"
$ws.Range("H135").Value = "750.00"
$ws.Range("H136").Value = "3500.00"
$ws.Range("H137").Value = "600.00"
$ws.Range("H138").Value = "8384.27"
$ws.Range("H139").Value = "537.00"
$ws.Range("H140").Value = "5445.00"
$ws.Range("H141").Value = "360.00"
$ws.Range("H142").Value = "360.00"
$ws.Range("H143").Value = "300.00"
$ws.Range("H144").Value = "3032.00"
$ws.Range("H145").Value = "12678.65"
$ws.Range("H146").Value = "1500.00"
$ws.Range("H147").Value = "650.00"
$ws.Range("H148").Value = "700.00"
$ws.Range("H149").Value = "120.00"
$ws.Range("H150").Value = "10929.46"
$ws.Range("H151").Value = "108.90"
$ws.Range("H152").Value = "107.00"
$ws.Range("H153").Value = "289.25"
$ws.Range("H154").Value = "130.00"
$ws.Range("H155").Value = "180.00"
$ws.Range("H156").Value = "162.00"
$ws.Range("H157").Value = "399.00"
$ws.Range("H158").Value = "260.00"
$ws.Range("H159").Value = "145.76"
$ws.Range("H160").Value = "2725.70"
$ws.Range("H161").Value = "227.06"
$ws.Range("H162").Value = "1008.00"
$ws.Range("H163").Value = "2580.33"
$ws.Range("H164").Value = "55.08"
$ws.Range("H165").Value = "1882.00"
$ws.Range("H166").Value = "8874.00"
$ws.Range("H167").Value = "1010.05"
$ws.Range("H168").Value = "180.00"
$ws.Range("H169").Value = "97.13"
$ws.Range("H170").Value = "5.00"
$ws.Range("H171").Value = "400.00"
$ws.Range("H172").Value = "384.90"
$ws.Range("H173").Value = "1809.00"
$ws.Range("H174").Value = "1675.00"
$ws.Range("H175").Value = "1917.00"
$ws.Range("H176").Value = "3700.79"
$ws.Range("H177").Value = "52532.61"
$ws.Range("H178").Value = "4950.00"
$ws.Range("H179").Value = "68911.91"
$ws.Range("H180").Value = "1415.00"
$ws.Range("H181").Value = "40500.00"
$ws.Range("H182").Value = "942.26"
$ws.Range("H183").Value = "16500.00"
$ws.Range("H184").Value = "3300.00"
$ws.Range("H185").Value = "2214.55"
$ws.Range("H186").Value = "22542.00"
$ws.Range("H187").Value = "364.00"
$ws.Range("H188").Value = "5000.00"
$ws.Range("H189").Value = "481.00"
$importeRange.Style = "Normal"
